$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add row 39: "39. String + DP" / Delete Operation for Two Strings ---
# Seed D39 by copying D38 (preserves the existing "Link" cell style s="9"
# without Excel synthesizing a brand-new cellXf for the hyperlink font),
# then repoint the copied hyperlink to the new URL.
$ws.Range("D38").Copy($ws.Range("D39"))

$link = $ws.Range("D39").Hyperlinks.Item(1)
$link.Address = "https://leetcode.com/problems/delete-operation-for-two-strings/"

$ws.Range("D39").Value = "https://leetcode.com/problems/delete-operation-for-two-strings/"
$ws.Range("A39").Value = "39. String + DP"
$ws.Range("B39").Value = "Delete Operation for Two Strings"
$ws.Range("C39").Value = "The minimum number of steps = length of word1 + length of word2 - (2 * length of LCS)"

# --- Update the saved view state (scrolled down a bit, new selection) ---
$ws.Range("B46").Select() | Out-Null
